$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "test"
$ws.Range("B5").Value = 400
$ws.Range("D5").Value = 200

$ws.Range("I10").Select()
